$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = "Dr. Servinaz Sayed Mohammad, Administrator, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy"
$ws.Range("G3").Value2 = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Veronia Rafat, Administrator"
$ws.Range("G4").Value2 = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy, Dr. Gehan Adel"
$ws.Range("G5").Value2 = "Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi"
$ws.Range("G6").Value2 = "Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany"
$ws.Range("G7").Value2 = "Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Kerelos Zareef, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G9").Value2 = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G11").Value2 = "Dr. Aya Saeed, Dr. Safa Hany, Dr. Amal Awwad"
$ws.Range("G12").Value2 = "Dr. Marina Youhanna, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G13").Value2 = "Dr. Esraa Mostafa, Dr. Amira Ibrahim, Dr. Yasmeena Fattoh"
$ws.Range("G20").Value2 = "Dr. Mohammad Safwat, Dr. Mariam Toma Gerges"
$ws.Range("G25").Value2 = "Menna tuâ€™Allah Gamil, Dr. Nouran Mahmoud"
$ws.Range("G30").Value2 = "Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Wafaa Ebida, Dr. Shorok Mohammad"
